$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 254 (shifts the existing rows 254-276 down to 255-277,
# carrying their values/styles with them).
$ws.Rows(254).Insert()

# Populate the newly inserted row 254 with the new weekly data point.
$ws.Range("A254").Value = 3
$ws.Range("B254").Value = "Femacal de La Calera"
$ws.Range("C254").Value = "Coquimbo"
$ws.Range("D254").Value = 44578
$ws.Range("E254").Value = 5
$ws.Range("F254").Value = 100114013
$ws.Range("G254").Value = "Zanahoria"
$ws.Range("H254").Value = "Sin especificar"
$ws.Range("I254").Value = "Primera"
$ws.Range("J254").Value = 530
$ws.Range("K254").Value = 7000
$ws.Range("L254").Value = 8000
$ws.Range("M254").Value = 7528
$ws.Range("N254").Value = "`$/saco 20 kilos"
$ws.Range("O254").Value = "Provincia de Quillota"
$ws.Range("P254").Value = 376
$ws.Range("Q254").Value = 20
$ws.Range("R254").Value = "Hortaliza"
